$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.00", "582.31") are not coerced into numbers, matching the
# original inline-string cell type used throughout the price column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 2).Value2 = 'Bitcoin'
$ws.Cells.Item(2, 3).Value2 = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Cells.Item(2, 4).Value2 = '66.968.53'
$ws.Cells.Item(2, 5).Value2 = '  +0.55%  '

$ws.Cells.Item(3, 2).Value2 = 'Ethereum'
$ws.Cells.Item(3, 3).Value2 = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Cells.Item(3, 4).Value2 = '2.472.71'
$ws.Cells.Item(3, 5).Value2 = '  +1.41%  '

$ws.Cells.Item(4, 2).Value2 = 'TetherUSD'
$ws.Cells.Item(4, 3).Value2 = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Cells.Item(4, 4).Value2 = '1.00'
$ws.Cells.Item(4, 5).Value2 = '  -0.04%  '

$ws.Cells.Item(5, 2).Value2 = 'BNB'
$ws.Cells.Item(5, 3).Value2 = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(5, 4).Value2 = '582.31'
$ws.Cells.Item(5, 5).Value2 = '  +0.98%  '

$ws.Cells.Item(6, 2).Value2 = 'Solana'
$ws.Cells.Item(6, 3).Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(6, 4).Value2 = '170.72'
$ws.Cells.Item(6, 5).Value2 = '  +3.73%  '

$ws.Cells.Item(7, 2).Value2 = 'USDC'
$ws.Cells.Item(7, 3).Value2 = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(7, 4).Value2 = '1.00'
$ws.Cells.Item(7, 5).Value2 = '  -0.10%  '

$ws.Cells.Item(8, 2).Value2 = 'XRP'
$ws.Cells.Item(8, 3).Value2 = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(8, 4).Value2 = '0.513'
$ws.Cells.Item(8, 5).Value2 = '  +1.19%  '

$ws.Cells.Item(9, 2).Value2 = 'LidoStakedEther'
$ws.Cells.Item(9, 3).Value2 = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(9, 4).Value2 = '2.474.05'
$ws.Cells.Item(9, 5).Value2 = '  +1.37%  '

$ws.Cells.Item(10, 2).Value2 = 'Dogecoin'
$ws.Cells.Item(10, 3).Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 4).Value2 = '0.136'
$ws.Cells.Item(10, 5).Value2 = '  +3.00%  '

$ws.Cells.Item(11, 2).Value2 = 'TRON'
$ws.Cells.Item(11, 3).Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(11, 4).Value2 = '0.166'
$ws.Cells.Item(11, 5).Value2 = '  +1.63%  '

$ws.Cells.Item(12, 2).Value2 = 'Toncoin'
$ws.Cells.Item(12, 3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(12, 4).Value2 = '4.91'
$ws.Cells.Item(12, 5).Value2 = '  +2.16%  '

$ws.Cells.Item(13, 2).Value2 = 'Cardano'
$ws.Cells.Item(13, 3).Value2 = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(13, 4).Value2 = '0.331'
$ws.Cells.Item(13, 5).Value2 = '  +1.41%  '

$ws.Cells.Item(14, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value2 = '2.976.37'
$ws.Cells.Item(14, 5).Value2 = '  +4.05%  '

$ws.Cells.Item(15, 2).Value2 = 'Avalanche'
$ws.Cells.Item(15, 3).Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value2 = '25.28'
$ws.Cells.Item(15, 5).Value2 = '  +1.19%  '

$ws.Cells.Item(16, 2).Value2 = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value2 = '66.897.53'
$ws.Cells.Item(16, 5).Value2 = '  +0.75%  '

$ws.Cells.Item(17, 2).Value2 = 'ShibaInu'
$ws.Cells.Item(17, 3).Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).Value2 = '0.0000169'
$ws.Cells.Item(17, 5).Value2 = '  +1.40%  '

$ws.Cells.Item(18, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(18, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value2 = '2.417.88'
$ws.Cells.Item(18, 5).Value2 = '  -0.97%  '

$ws.Cells.Item(19, 2).Value2 = 'Chainlink'
$ws.Cells.Item(19, 3).Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(19, 4).Value2 = '10.91'
$ws.Cells.Item(19, 5).Value2 = '  -2.47%  '

$ws.Cells.Item(20, 2).Value2 = 'Uniswap'
$ws.Cells.Item(20, 3).Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).Value2 = '7.41'
$ws.Cells.Item(20, 5).Value2 = '  -0.73%  '

$ws.Cells.Item(21, 2).Value2 = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(21, 4).Value2 = '348.84'
$ws.Cells.Item(21, 5).Value2 = '  -0.63%  '

$ws.Cells.Item(22, 2).Value2 = 'Polkadot'
$ws.Cells.Item(22, 3).Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(22, 4).Value2 = '4.02'
$ws.Cells.Item(22, 5).Value2 = '  +0.73%  '

$ws.Cells.Item(23, 2).Value2 = 'Dai'
$ws.Cells.Item(23, 3).Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value2 = '1.00'
$ws.Cells.Item(23, 5).Value2 = '  -0.16%  '

$ws.Cells.Item(24, 2).Value2 = 'Litecoin'
$ws.Cells.Item(24, 3).Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).Value2 = '68.15'
$ws.Cells.Item(24, 5).Value2 = '  -1.68%  '

$ws.Cells.Item(25, 2).Value2 = 'NEARProtocol'
$ws.Cells.Item(25, 3).Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(25, 4).Value2 = '4.21'
$ws.Cells.Item(25, 5).Value2 = '  +0.60%  '

$ws.Cells.Item(26, 2).Value2 = 'SuiNetwork'
$ws.Cells.Item(26, 3).Value2 = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(26, 4).Value2 = '1.78'
$ws.Cells.Item(26, 5).Value2 = '  +3.22%  '

$ws.Cells.Item(27, 2).Value2 = 'Aptos'
$ws.Cells.Item(27, 3).Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(27, 4).Value2 = '9.28'
$ws.Cells.Item(27, 5).Value2 = '  +4.50%  '

$ws.Cells.Item(28, 2).Value2 = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value2 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value2 = '2.606.89'
$ws.Cells.Item(28, 5).Value2 = '  +1.57%  '

$ws.Cells.Item(29, 2).Value2 = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value2 = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(29, 4).Value2 = '0.998'
$ws.Cells.Item(29, 5).Value2 = '  -0.08%  '

$ws.Cells.Item(30, 2).Value2 = 'PEPE'
$ws.Cells.Item(30, 3).Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(30, 4).Value2 = '0.0₃0898'
$ws.Cells.Item(30, 5).Value2 = '  +1.68%  '

$ws.Cells.Item(31, 2).Value2 = 'Bittensor'
$ws.Cells.Item(31, 3).Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).Value2 = '505.79'
$ws.Cells.Item(31, 5).Value2 = '  +0.31%  '

$ws.Cells.Item(32, 2).Value2 = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value2 = '7.69'
$ws.Cells.Item(32, 5).Value2 = '  -0.22%  '

$ws.Cells.Item(33, 2).Value2 = 'Fetch.AI'
$ws.Cells.Item(33, 3).Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(33, 4).Value2 = '1.23'
$ws.Cells.Item(33, 5).Value2 = '  +1.79%  '

$ws.Cells.Item(34, 2).Value2 = 'PancakeSwap'
$ws.Cells.Item(34, 3).Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(34, 4).Value2 = '1.75'
$ws.Cells.Item(34, 5).Value2 = '  +0.03%  '

$ws.Cells.Item(35, 2).Value2 = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).Value2 = '1.00'
$ws.Cells.Item(35, 5).Value2 = '  -0.02%  '

$ws.Cells.Item(36, 2).Value2 = 'Monero'
$ws.Cells.Item(36, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).Value2 = '160.58'
$ws.Cells.Item(36, 5).Value2 = '  +0.72%  '

$ws.Cells.Item(37, 2).Value2 = 'Kaspa'
$ws.Cells.Item(37, 3).Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).Value2 = '0.117'
$ws.Cells.Item(37, 5).Value2 = '  +3.29%  '

$ws.Cells.Item(38, 2).Value2 = 'WhiteBITCoin'
$ws.Cells.Item(38, 3).Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(38, 4).Value2 = '18.69'
$ws.Cells.Item(38, 5).Value2 = '  +0.83%  '

$ws.Cells.Item(39, 2).Value2 = 'EthereumClassic'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(39, 4).Value2 = '18.19'
$ws.Cells.Item(39, 5).Value2 = '  -0.73%  '

$ws.Cells.Item(40, 2).Value2 = 'ImmutableX'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(40, 4).Value2 = '1.33'
$ws.Cells.Item(40, 5).Value2 = '  -0.30%  '

$ws.Cells.Item(41, 2).Value2 = 'USDe'
$ws.Cells.Item(41, 3).Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(41, 4).Value2 = '1.00'
$ws.Cells.Item(41, 5).Value2 = '  +0.03%  '

$ws.Cells.Item(42, 2).Value2 = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value2 = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(42, 4).Value2 = '0.328'
$ws.Cells.Item(42, 5).Value2 = '  +1.64%  '

$ws.Cells.Item(43, 2).Value2 = 'Stacks'
$ws.Cells.Item(43, 3).Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).Value2 = '1.68'
$ws.Cells.Item(43, 5).Value2 = '  +0.99%  '

$ws.Cells.Item(44, 2).Value2 = 'RenderToken'
$ws.Cells.Item(44, 3).Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(44, 4).Value2 = '4.78'
$ws.Cells.Item(44, 5).Value2 = '  +1.89%  '

$ws.Cells.Item(45, 2).Value2 = 'dogwifhat'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(45, 4).Value2 = '2.35'
$ws.Cells.Item(45, 5).Value2 = '  +2.37%  '

$ws.Cells.Item(46, 2).Value2 = 'Aave'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).Value2 = '142.82'
$ws.Cells.Item(46, 5).Value2 = '  +1.90%  '

$ws.Cells.Item(47, 2).Value2 = 'Filecoin'
$ws.Cells.Item(47, 3).Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(47, 4).Value2 = '3.46'
$ws.Cells.Item(47, 5).Value2 = '  +0.48%  '

$ws.Cells.Item(48, 2).Value2 = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).Value2 = '0.512'
$ws.Cells.Item(48, 5).Value2 = '  +1.13%  '

$ws.Cells.Item(49, 2).Value2 = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(49, 4).Value2 = '0.0₆0250'
$ws.Cells.Item(49, 5).Value2 = '  +4.25%  '

$ws.Cells.Item(50, 2).Value2 = 'Cronos'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value2 = '0.0730'
$ws.Cells.Item(50, 5).Value2 = '  +0.31%  '

$ws.Cells.Item(51, 2).Value2 = 'Optimism'
$ws.Cells.Item(51, 3).Value2 = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Cells.Item(51, 4).Value2 = '1.57'
$ws.Cells.Item(51, 5).Value2 = '  +0.00%  '
